$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1169995834814548
$ws.Range("C2").Value = 41249014.21622031
$ws.Range("D2").Value = 18.71679738969934
$ws.Range("E2").Value = 14773364.14517103
$ws.Range("G2").Value = 56022397.19518831

$ws.Range("B3").Value = 3.272327238179451
$ws.Range("C3").Value = 0.3048912486333797
$ws.Range("D3").Value = 0.7210945179870265
$ws.Range("E3").Value = 2797.565817734744
$ws.Range("G3").Value = 2801.864130739543

$ws.Range("B4").Value = 0.1169995834814548
$ws.Range("C4").Value = 9.983522426115931
$ws.Range("D4").Value = 0.7210945179870265
$ws.Range("E4").Value = 2797.565817734744
$ws.Range("G4").Value = 2808.387434262328
